$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.146.49"
$ws.Range("E2").Value = "  -6.89%  "
$ws.Range("D3").Value = "2.552.19"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("D5").Value = "'298.77"
$ws.Range("E5").Value = "  -4.29%  "
$ws.Range("D6").Value = "'94.43"
$ws.Range("E6").Value = "  -5.63%  "
$ws.Range("E7").Value = "  -3.91%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.554"
$ws.Range("E9").Value = "  -5.67%  "
$ws.Range("D10").Value = "'36.18"
$ws.Range("E10").Value = "  -7.83%  "
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("D12").Value = "'7.77"
$ws.Range("E12").Value = "  -5.06%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "2.939.21"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "2.539.05"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "'0.872"
$ws.Range("E16").Value = "  -5.57%  "
$ws.Range("D17").Value = "'14.23"
$ws.Range("E17").Value = "  -5.03%  "
$ws.Range("D18").Value = "43.125.23"
$ws.Range("E18").Value = "  -7.61%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.68"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "'12.69"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "'71.89"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'261.95"
$ws.Range("E23").Value = "  -10.82%  "
$ws.Range("E24").Value = "  -5.24%  "
$ws.Range("B25").Value = "EthereumClassic"
$ws.Range("C25").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D25").Value = "'29.72"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'2.14"
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("E28").Value = "  -7.45%  "
$ws.Range("D29").Value = "'37.00"
$ws.Range("E29").Value = "  -6.44%  "
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").Value = "'5.98"
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("D32").Value = "'155.33"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("D35").Value = "'2.74"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "'0.0801"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("E37").Value = "  -5.80%  "
$ws.Range("D38").Value = "'0.120"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("D39").Value = "'23.55"
$ws.Range("E39").Value = "  +9.96%  "
$ws.Range("D40").Value = "'16.59"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("E42").Value = "  -6.07%  "
$ws.Range("D43").Value = "'3.88"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").Value = "2.084.43"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'86.17"
$ws.Range("E46").Value = "  -12.09%  "
$ws.Range("D47").Value = "'1.59"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("D48").Value = "2.794.66"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.70"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'104.65"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = "  -9.45%  "
